$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

# Row 2: D2->62.855.79
$ws.Range("D2").Value = '62.855.79'
$ws.Range("E2").Value = '  +4.24%  '

# Row 3: D3->2.461.42
$ws.Range("D3").Value = '2.461.42'
$ws.Range("E3").Value = '  +5.35%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5: D5->563.52
Set-TextCell "D5" '563.52'
$ws.Range("E5").Value = '  +2.62%  '

# Row 6: D6->141.85
Set-TextCell "D6" '141.85'
$ws.Range("E6").Value = '  +8.07%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8: D8->0.587
Set-TextCell "D8" '0.587'
$ws.Range("E8").Value = '  +1.17%  '

# Row 9: D9->2.460.00
$ws.Range("D9").Value = '2.460.00'
$ws.Range("E9").Value = '  +5.43%  '

# Row 10: D10->0.105
Set-TextCell "D10" '0.105'
$ws.Range("E10").Value = '  +2.83%  '

# Row 11: D11->5.69
Set-TextCell "D11" '5.69'
$ws.Range("E11").Value = '  +1.48%  '

# Row 12
$ws.Range("E12").Value = '  +1.04%  '

# Row 13
$ws.Range("E13").Value = '  +4.38%  '

# Row 14: D14->26.21
Set-TextCell "D14" '26.21'
$ws.Range("E14").Value = '  +10.26%  '

# Row 15: D15->2.903.62
$ws.Range("D15").Value = '2.903.62'
$ws.Range("E15").Value = '  +5.58%  '

# Row 16: D16->62.809.12
$ws.Range("D16").Value = '62.809.12'
$ws.Range("E16").Value = '  +4.29%  '

# Row 17: D17->0.0000140
Set-TextCell "D17" '0.0000140'
$ws.Range("E17").Value = '  +4.24%  '

# Row 18: D18->2.464.48
$ws.Range("D18").Value = '2.464.48'
$ws.Range("E18").Value = '  +5.76%  '

# Row 19
$ws.Range("E19").Value = '  +5.45%  '

# Row 20: D20->339.30
Set-TextCell "D20" '339.30'
$ws.Range("E20").Value = '  +7.98%  '

# Row 21: D21->4.25
Set-TextCell "D21" '4.25'
$ws.Range("E21").Value = '  +3.53%  '

# Row 22: D22->6.78
Set-TextCell "D22" '6.78'
$ws.Range("E22").Value = '  +2.77%  '

# Row 23: D23->0.998
Set-TextCell "D23" '0.998'
$ws.Range("E23").Value = '  -0.15%  '

# Row 24: D24->65.40
Set-TextCell "D24" '65.40'
$ws.Range("E24").Value = '  +2.11%  '

# Row 25: D25->0.174
Set-TextCell "D25" '0.174'
$ws.Range("E25").Value = '  +1.88%  '

# Row 26
$ws.Range("E26").Value = '  +0.07%  '

# Row 27: D27->1.49
Set-TextCell "D27" '1.49'
$ws.Range("E27").Value = '  +7.78%  '

# Row 28: D28->8.04
Set-TextCell "D28" '8.04'
$ws.Range("E28").Value = '  +1.32%  '

# Row 29: D29->1.36
Set-TextCell "D29" '1.36'
$ws.Range("E29").Value = '  +8.10%  '

# Row 30: D30->6.80
Set-TextCell "D30" '6.80'
$ws.Range("E30").Value = '  +11.80%  '

# Row 31: B31->PancakeSwap, D31->1.84
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell "D31" '1.84'
$ws.Range("E31").Value = '  +6.13%  '

# Row 32: B32->PEPE, D32->0.0₃0797
$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").Value = '0.0₃0797'
$ws.Range("E32").Value = '  +8.65%  '

# Row 33: D33->176.92
Set-TextCell "D33" '176.92'
$ws.Range("E33").Value = '  +3.31%  '

# Row 34
$ws.Range("E34").Value = '  +10.32%  '

# Row 35: D35->0.395
Set-TextCell "D35" '0.395'
$ws.Range("E35").Value = '  +2.84%  '

# Row 36: D36->18.79
Set-TextCell "D36" '18.79'
$ws.Range("E36").Value = '  +3.94%  '

# Row 37: D37->363.73
Set-TextCell "D37" '363.73'
$ws.Range("E37").Value = '  +12.07%  '

# Row 38: B38->USDe, D38->0.999
$ws.Range("B38").Value = 'USDe'
$ws.Range("C38").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell "D38" '0.999'
$ws.Range("E38").Value = '  +0.01%  '

# Row 39: B39->NEARProtocol, D39->4.37
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell "D39" '4.37'
$ws.Range("E39").Value = '  +6.32%  '

# Row 40: D40->0.999
Set-TextCell "D40" '0.999'
$ws.Range("E40").Value = '  -0.03%  '

# Row 41
$ws.Range("E41").Value = '  +10.30%  '

# Row 42: D42->40.51
Set-TextCell "D42" '40.51'
$ws.Range("E42").Value = '  +6.13%  '

# Row 43: D43->149.15
Set-TextCell "D43" '149.15'
$ws.Range("E43").Value = '  +8.29%  '

# Row 44: D44->3.69
Set-TextCell "D44" '3.69'
$ws.Range("E44").Value = '  +5.30%  '

# Row 45: D45->20.50
Set-TextCell "D45" '20.50'
$ws.Range("E45").Value = '  +6.00%  '

# Row 46: D46->0.595
Set-TextCell "D46" '0.595'
$ws.Range("E46").Value = '  +4.88%  '

# Row 47: D47->0.0955
Set-TextCell "D47" '0.0955'
$ws.Range("E47").Value = '  +0.42%  '

# Row 48: D48->0.0513
Set-TextCell "D48" '0.0513'
$ws.Range("E48").Value = '  +3.20%  '

# Row 49: D49->0.0₆0236
$ws.Range("D49").Value = '0.0₆0236'
$ws.Range("E49").Value = '  +7.56%  '

# Row 50: D50->0.0225
Set-TextCell "D50" '0.0225'
$ws.Range("E50").Value = '  +4.12%  '

# Row 51: D51->17.87
Set-TextCell "D51" '17.87'
$ws.Range("E51").Value = '  +4.56%  '
